$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column A (rankings list got longer / url prefix changed)
$ws.Columns.Item(1).ColumnWidth = 37.83

# Re-order the QB ranking list (rows 3-16) to reflect the new ranking order
$ws.Range("A3").Value = "Drew Brees - NO"
$ws.Range("A4").Value = "Russell Wilson - SEA"
$ws.Range("A5").Value = "Andrew Luck - IND"
$ws.Range("A6").Value = "Tom Brady - NE"
$ws.Range("A7").Value = "Philip Rivers - LAC"
$ws.Range("A8").Value = "Ben Roethlisberger - PIT"
$ws.Range("A9").Value = "Eli Manning - NYG"
$ws.Range("A10").Value = "Matthew Stafford - DET"
$ws.Range("A11").Value = "Kirk Cousins - MIN"
$ws.Range("A12").Value = "Deshaun Watson - HOU"
$ws.Range("A13").Value = "Matt Ryan - ATL"
$ws.Range("A14").Value = "Jared Goff - LAR"
$ws.Range("A15").Value = "Marcus Mariota - TEN"
$ws.Range("A16").Value = "Carson Wentz - PHI"
